$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# Insert a new data row for "VONACIDAN 20 MG 20 F.C.TABS." right
# before "WATER FOR INJECTION AMP. 5 ML" (row 107), keeping the
# shortage list in alphabetical order.
# -----------------------------------------------------------------

$ws.Rows("107:107").Insert()

# Clone the formatting (styles, number formats, fonts, ...) of the
# row above so the new row matches the rest of the table exactly.
$ws.Range("A106:Q106").Copy()
$ws.Range("A107:Q107").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Match row height used by the table rows.
$ws.Rows(107).RowHeight = $ws.Rows(106).RowHeight

# Recreate the merged cells for the new row (A:B, C:G, H:K, L:M, N:O)
$ws.Range("A107:B107").Merge()
$ws.Range("C107:G107").Merge()
$ws.Range("H107:K107").Merge()
$ws.Range("L107:M107").Merge()
$ws.Range("N107:O107").Merge()

function Set-TextValue($cell, $text) {
    $fmt = $cell.NumberFormat
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.NumberFormat = $fmt
}

# Item name
Set-TextValue $ws.Cells.Item(107, 3) "VONACIDAN 20 MG 20 F.C.TABS."
# Current balance ratio
Set-TextValue $ws.Cells.Item(107, 8) "0:0"
# Price
Set-TextValue $ws.Cells.Item(107, 14) "138.00"
# Selling price
Set-TextValue $ws.Cells.Item(107, 16) "69.0000"
# Number-of-transactions ratio
Set-TextValue $ws.Cells.Item(107, 17) "0:1"

# -----------------------------------------------------------------
# Column A holds the row's running serial number (row - 6). Inserting
# the row shifted the previously-stored numbers down by one row along
# with their cells, so every serial number from the new row through
# the last data row needs to be restored back to "row - 6".
# -----------------------------------------------------------------
for ($r = 107; $r -le 131; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 6
}

# -----------------------------------------------------------------
# Update the grand-total cell (now row 132) to include the new item's
# selling price, and refresh the generated-on timestamp (now row 133).
# -----------------------------------------------------------------
$ws.Cells.Item(132, 16).Value = $ws.Cells.Item(132, 16).Value2 + 69
$ws.Cells.Item(133, 1).Value = "Friday, 22 August, 2025 10:56 PM"

$ws.Application.CutCopyMode = $false
